$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds textual price figures (e.g. "1.00", "45.086.31") that Excel would
# otherwise coerce into numbers. Temporarily force text format while writing those
# values, then clear the formatting again so cells keep their original (unset) style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '45.086.31'
$ws.Range("D3").Value = '2.381.47'
$ws.Range("D5").Value = '293.05'
$ws.Range("D6").Value = '94.34'
$ws.Range("D7").Value = '0.554'
$ws.Range("D9").Value = '0.496'
$ws.Range("D10").Value = '33.92'
$ws.Range("D11").Value = '0.0773'
$ws.Range("D12").Value = '6.92'
$ws.Range("D14").Value = '2.744.96'
$ws.Range("D15").Value = '2.391.13'
$ws.Range("D16").Value = '13.90'
$ws.Range("D17").Value = '0.819'
$ws.Range("D18").Value = '45.052.10'
$ws.Range("D19").Value = '12.32'
$ws.Range("D20").Value = '0.0₃0923'
$ws.Range("D22").Value = '66.07'
$ws.Range("D23").Value = '237.68'
$ws.Range("D25").Value = '1.00'
$ws.Range("D27").Value = '2.19'
$ws.Range("D28").Value = '37.24'
$ws.Range("D29").Value = '9.48'
$ws.Range("D31").Value = '20.90'
$ws.Range("D32").Value = '2.71'
$ws.Range("D33").Value = '146.99'
$ws.Range("D34").Value = '5.37'
$ws.Range("D35").Value = '0.0754'
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.112'
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = '1.94'
$ws.Range("D39").Value = '14.67'
$ws.Range("D40").Value = '3.71'
$ws.Range("D41").Value = '0.0291'
$ws.Range("D42").Value = '3.15'
$ws.Range("D43").Value = '1.942.99'
$ws.Range("D44").Value = '0.999'
$ws.Range("D45").Value = '88.88'
$ws.Range("D46").Value = '1.72'
$ws.Range("D47").Value = '8.42'
$ws.Range("D48").Value = '14.97'
$ws.Range("D49").Value = '98.97'
$ws.Range("D50").Value = '2.615.36'
$ws.Range("D51").Value = '0.181'

$ws.Range("D2:D51").ClearFormats()

# Column E (percentage change) never parses as a plain number because of the "%"
# sign and surrounding spaces, so it can be written directly.
$ws.Range("E2").Value = '  -3.35%  '
$ws.Range("E3").Value = '  +4.99%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E5").Value = '  -3.20%  '
$ws.Range("E6").Value = '  -5.69%  '
$ws.Range("E7").Value = '  -1.38%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -3.41%  '
$ws.Range("E10").Value = '  -4.95%  '
$ws.Range("E11").Value = '  -1.31%  '
$ws.Range("E12").Value = '  -3.91%  '
$ws.Range("E13").Value = '  +0.75%  '
$ws.Range("E14").Value = '  +4.87%  '
$ws.Range("E15").Value = '  +4.93%  '
$ws.Range("E16").Value = '  +1.64%  '
$ws.Range("E17").Value = '  +2.35%  '
$ws.Range("E18").Value = '  -3.41%  '
$ws.Range("E19").Value = '  -6.46%  '
$ws.Range("E20").Value = '  -0.65%  '
$ws.Range("E21").Value = '  +2.54%  '
$ws.Range("E22").Value = '  +1.21%  '
$ws.Range("E23").Value = '  -3.90%  '
$ws.Range("E24").Value = '  -4.02%  '
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("E28").Value = '  -12.87%  '
$ws.Range("E29").Value = '  -2.61%  '
$ws.Range("E30").Value = '  +15.85%  '
$ws.Range("E31").Value = '  +5.22%  '
$ws.Range("E32").Value = '  -2.50%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("E34").Value = '  -2.41%  '
$ws.Range("E35").Value = '  -2.90%  '
$ws.Range("E36").Value = '  -3.35%  '
$ws.Range("E37").Value = '  +12.20%  '
$ws.Range("E38").Value = '  -2.07%  '
$ws.Range("E39").Value = '  -8.36%  '
$ws.Range("E40").Value = '  -5.21%  '
$ws.Range("E41").Value = '  -1.85%  '
$ws.Range("E42").Value = '  -3.35%  '
$ws.Range("E43").Value = '  +7.05%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("E46").Value = '  -13.18%  '
$ws.Range("E47").Value = '  +7.68%  '
$ws.Range("E48").Value = '  +15.64%  '
$ws.Range("E49").Value = '  +4.77%  '
$ws.Range("E50").Value = '  +4.84%  '
$ws.Range("E51").Value = '  -4.65%  '
